# A new weekly price observation ("Chino" / "Primera", 2023-12-dated) was
# inserted into the daily log as row 285, pushing every existing row
# (285-355) down by one (286-356).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 285; rows 285:355 shift down to 286:356,
# carrying their formatting (incl. the date-styled column D) with them.
$ws.Rows("285:285").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A285").Value = 9
$ws.Range("B285").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C285").Value = "Metropolitana"
$ws.Range("D285").Value = 45173
$ws.Range("E285").Value = 13
$ws.Range("F285").Value = 100112003
$ws.Range("G285").Value = "Ajo"
$ws.Range("H285").Value = "Chino"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 610
$ws.Range("K285").Value = 17000
$ws.Range("L285").Value = 18000
$ws.Range("M285").Value = 17500
$ws.Range("N285").Value = "`$/caja 10 kilos"
$ws.Range("O285").Value = "China"
$ws.Range("P285").Value = 1750
$ws.Range("Q285").Value = 10
$ws.Range("R285").Value = "Hortaliza"
